$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.314.55'
$ws.Range('E2').Value = '  -1.33%  '

# Row 3
$ws.Range('D3').Value = '3.004.94'
$ws.Range('E3').Value = '  -1.42%  '

# Row 4
$ws.Range('E4').Value = '  +0.21%  '

# Row 5
$ws.Range('D5').Value = '586.51'
$ws.Range('E5').Value = '  -0.24%  '

# Row 6
$ws.Range('D6').Value = '146.39'
$ws.Range('E6').Value = '  -3.19%  '

# Row 7
$ws.Range('E7').Value = '  -0.04%  '

# Row 8
$ws.Range('D8').Value = '0.526'
$ws.Range('E8').Value = '  -2.28%  '

# Row 9
$ws.Range('D9').Value = '3.005.61'
$ws.Range('E9').Value = '  -1.33%  '

# Row 10
$ws.Range('D10').Value = '0.148'
$ws.Range('E10').Value = '  -3.76%  '

# Row 11
$ws.Range('D11').Value = '5.78'
$ws.Range('E11').Value = '  -0.16%  '

# Row 12
$ws.Range('D12').Value = '0.463'
$ws.Range('E12').Value = '  +3.09%  '

# Row 13
$ws.Range('D13').Value = '0.0000229'
$ws.Range('E13').Value = '  -2.35%  '

# Row 14
$ws.Range('D14').Value = '34.52'

# Row 15
$ws.Range('D15').Value = '0.123'
$ws.Range('E15').Value = '  +2.12%  '

# Row 16
$ws.Range('D16').Value = '3.514.68'
$ws.Range('E16').Value = '  -1.04%  '

# Row 17
$ws.Range('D17').Value = '7.06'
$ws.Range('E17').Value = '  -1.21%  '

# Row 18
$ws.Range('D18').Value = '62.383.14'
$ws.Range('E18').Value = '  -1.26%  '

# Row 19
$ws.Range('D19').Value = '3.011.64'
$ws.Range('E19').Value = '  -1.36%  '

# Row 20
$ws.Range('D20').Value = '458.07'
$ws.Range('E20').Value = '  -3.94%  '

# Row 21
$ws.Range('D21').Value = '13.98'
$ws.Range('E21').Value = '  -1.99%  '

# Row 22
$ws.Range('D22').Value = '0.688'
$ws.Range('E22').Value = '  -2.37%  '

# Row 23
$ws.Range('D23').Value = '7.40'
$ws.Range('E23').Value = '  -1.47%  '

# Row 24
$ws.Range('D24').Value = '81.63'
$ws.Range('E24').Value = '  -0.61%  '

# Row 25
$ws.Range('D25').Value = '2.20'
$ws.Range('E25').Value = '  -8.90%  '

# Row 26
$ws.Range('D26').Value = '12.21'
$ws.Range('E26').Value = '  -3.90%  '

# Row 27
$ws.Range('E27').Value = '  -0.08%  '

# Row 28
$ws.Range('D28').Value = '9.84'
$ws.Range('E28').Value = '  -7.79%  '

# Row 29
$ws.Range('E29').Value = '  +0.19%  '

# Row 30
$ws.Range('D30').Value = '2.64'
$ws.Range('E30').Value = '  -1.38%  '

# Row 31
$ws.Range('D31').Value = '6.96'
$ws.Range('E31').Value = '  -5.31%  '

# Row 32
$ws.Range('D32').Value = '2.09'
$ws.Range('E32').Value = '  -4.84%  '

# Row 33
$ws.Range('D33').Value = '27.69'
$ws.Range('E33').Value = '  +0.07%  '

# Row 34
$ws.Range('D34').Value = '0.109'
$ws.Range('E34').Value = '  -1.35%  '

# Row 35
$ws.Range('D35').Value = '0.0₃0807'
$ws.Range('E35').Value = '  -1.23%  '

# Row 36
$ws.Range('E36').Value = '  -2.76%  '

# Row 37
$ws.Range('D37').Value = '5.75'
$ws.Range('E37').Value = '  -2.48%  '

# Row 38
$ws.Range('D38').Value = '2.11'
$ws.Range('E38').Value = '  -5.23%  '

# Row 39
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '50.45'
$ws.Range('E39').Value = '  -0.12%  '

# Row 40
$ws.Range('B40').Value = 'Cosmos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D40').Value = '9.18'
$ws.Range('E40').Value = '  -0.76%  '

# Row 41
$ws.Range('D41').Value = '0.123'
$ws.Range('E41').Value = '  +7.80%  '

# Row 42
$ws.Range('D42').Value = '2.89'
$ws.Range('E42').Value = '  -10.66%  '

# Row 43
$ws.Range('D43').Value = '393.95'
$ws.Range('E43').Value = '  -9.35%  '

# Row 44
$ws.Range('D44').Value = '0.0357'
$ws.Range('E44').Value = '  -1.61%  '

# Row 45
$ws.Range('D45').Value = '0.267'
$ws.Range('E45').Value = '  -7.05%  '

# Row 46
$ws.Range('D46').Value = '2.734.24'
$ws.Range('E46').Value = '  -3.31%  '

# Row 47
$ws.Range('D47').Value = '37.39'
$ws.Range('E47').Value = '  -2.36%  '

# Row 48
$ws.Range('D48').Value = '129.10'
$ws.Range('E48').Value = '  -0.57%  '

# Row 50
$ws.Range('E50').Value = '  -0.47%  '

# Row 51
$ws.Range('D51').Value = '2.19'
$ws.Range('E51').Value = '  -0.97%  '
